$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: label | is | VIR SEPA EXPLEO FRANCE | category | Salaire Diego
$ws.Range("A2").Value = "label"
$ws.Range("B2").Value = "is"
$ws.Range("C2").Value = "VIR SEPA EXPLEO FRANCE"
$ws.Range("D2").Value = "category"
$ws.Range("E2").Value = "Salaire Diego"

# Row 3: label | is | VIR SEPA EXPLEO FRANCE | supplierFound | Expleo
$ws.Range("A3").Value = "label"
$ws.Range("B3").Value = "is"
$ws.Range("C3").Value = "VIR SEPA EXPLEO FRANCE"
$ws.Range("D3").Value = "supplierFound"
$ws.Range("E3").Value = "Expleo"

# Row 4 (new): label | contains | CSE EXPLEO | supplierFound | CSE Expleo
$ws.Range("A4").Value = "label"
$ws.Range("B4").Value = "contains"
$ws.Range("C4").Value = "CSE EXPLEO"
$ws.Range("D4").Value = "supplierFound"
$ws.Range("E4").Value = "CSE Expleo"

# Widen column C to fit the new, longer values
$ws.Columns.Item(3).ColumnWidth = 22

# Update the active selection to match the saved view
$ws.Range("D11").Select()
